$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Codebook")

# Update the "mktD" row's label/notes cell (C9) with the corrected logistic
# deterrence parameters (ef = .014891 for cartels and .025906 for mergers).
$ws.Range("C9").Value = 'Estimated size of the market deterred by the decision on case "k". Deterrent effects are calculated according to a logistic function based on the 64 sector aggregation. [Logistic parameters: chi = 100, ef = .014891 (cartels) and .025906 (mergers), ny = 1]'

# Move the active selection/view back up to C8, as in the final report pass.
$ws.Range("C8").Select()
